$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1993670886075949
$ws.Range("C2").Value = 0.5664556962025317
$ws.Range("J2").Value = 0.02531645569620253
$ws.Range("P2").Value = 0.1329113924050633
$ws.Range("S2").Value = 0.0759493670886076
$ws.Range("B3").Value = 0.01530612244897959
$ws.Range("C3").Value = 0.0663265306122449
$ws.Range("J3").Value = 0.06122448979591837
$ws.Range("P3").Value = 0.7040816326530612
$ws.Range("S3").Value = 0.1530612244897959
$ws.Range("J4").Value = 0.02777777777777778
$ws.Range("P4").Value = 0.75
$ws.Range("S4").Value = 0.2222222222222222
$ws.Range("B6").Value = 0.07727272727272727
$ws.Range("D6").Value = 0.01818181818181818
$ws.Range("E6").Value = 0.004545454545454545
$ws.Range("F6").Value = 0.06818181818181818
$ws.Range("J6").Value = 0.1818181818181818
$ws.Range("O6").Value = 0.05
$ws.Range("Q6").Value = 0.1409090909090909
$ws.Range("R6").Value = 0.07272727272727272
$ws.Range("S6").Value = 0.3863636363636364
$ws.Range("B7").Value = 0.09139784946236559
$ws.Range("D7").Value = 0.02150537634408602
$ws.Range("F7").Value = 0.05376344086021505
$ws.Range("J7").Value = 0.1451612903225807
$ws.Range("O7").Value = 0.02150537634408602
$ws.Range("Q7").Value = 0.1774193548387097
$ws.Range("R7").Value = 0.04301075268817205
$ws.Range("S7").Value = 0.4462365591397849
$ws.Range("B8").Value = 0.1149144254278729
$ws.Range("D8").Value = 0.019559902200489
$ws.Range("E8").Value = 0.002444987775061125
$ws.Range("F8").Value = 0.07823960880195599
$ws.Range("J8").Value = 0.1149144254278729
$ws.Range("O8").Value = 0.02444987775061125
$ws.Range("Q8").Value = 0.1784841075794621
$ws.Range("R8").Value = 0.08312958435207823
$ws.Range("S8").Value = 0.3838630806845966
$ws.Range("B9").Value = 0.1124260355029586
$ws.Range("D9").Value = 0.01183431952662722
$ws.Range("F9").Value = 0.1183431952662722
$ws.Range("J9").Value = 0.08284023668639054
$ws.Range("O9").Value = 0.01183431952662722
$ws.Range("Q9").Value = 0.1183431952662722
$ws.Range("R9").Value = 0.1242603550295858
$ws.Range("S9").Value = 0.4201183431952663
$ws.Range("B10").Value = 0.1145339652448657
$ws.Range("D10").Value = 0.01579778830963665
$ws.Range("E10").Value = 0.00315955766192733
$ws.Range("F10").Value = 0.06556082148499211
$ws.Range("J10").Value = 0.1303317535545024
$ws.Range("O10").Value = 0.02132701421800948
$ws.Range("Q10").Value = 0.1951026856240126
$ws.Range("R10").Value = 0.08846761453396525
$ws.Range("S10").Value = 0.3657187993680885
$ws.Range("G11").Value = 0.1640866873065016
$ws.Range("J11").Value = 0.1145510835913313
$ws.Range("K11").Value = 0.2414860681114551
$ws.Range("L11").Value = 0.4705882352941176
$ws.Range("S11").Value = 0.009287925696594427
$ws.Range("G12").Value = 0.6909090909090909
$ws.Range("J12").Value = 0.1878787878787879
$ws.Range("K12").Value = 0.03636363636363636
$ws.Range("L12").Value = 0.04848484848484848
$ws.Range("S12").Value = 0.03636363636363636
$ws.Range("G13").Value = 0.7586206896551724
$ws.Range("J13").Value = 0.2413793103448276
$ws.Range("F15").Value = 0.01345291479820628
$ws.Range("H15").Value = 0.08968609865470852
$ws.Range("I15").Value = 0.07174887892376682
$ws.Range("J15").Value = 0.3901345291479821
$ws.Range("K15").Value = 0.06278026905829596
$ws.Range("M15").Value = 0.0179372197309417
$ws.Range("O15").Value = 0.08520179372197309
$ws.Range("S15").Value = 0.2690582959641256
$ws.Range("F16").Value = 0.01492537313432836
$ws.Range("H16").Value = 0.154228855721393
$ws.Range("I16").Value = 0.07960199004975124
$ws.Range("J16").Value = 0.4129353233830846
$ws.Range("K16").Value = 0.1243781094527363
$ws.Range("M16").Value = 0.01990049751243781
$ws.Range("O16").Value = 0.06965174129353234
$ws.Range("S16").Value = 0.1243781094527363
$ws.Range("F17").Value = 0.01985111662531017
$ws.Range("H17").Value = 0.1935483870967742
$ws.Range("I17").Value = 0.09925558312655088
$ws.Range("J17").Value = 0.3895781637717122
$ws.Range("K17").Value = 0.08436724565756824
$ws.Range("M17").Value = 0.007444168734491315
$ws.Range("N17").Value = 0.004962779156327543
$ws.Range("O17").Value = 0.04962779156327544
$ws.Range("S17").Value = 0.1513647642679901
$ws.Range("F18").Value = 0.03684210526315789
$ws.Range("H18").Value = 0.1947368421052632
$ws.Range("I18").Value = 0.05263157894736842
$ws.Range("J18").Value = 0.4210526315789473
$ws.Range("K18").Value = 0.09473684210526316
$ws.Range("M18").Value = 0.01052631578947368
$ws.Range("O18").Value = 0.05263157894736842
$ws.Range("S18").Value = 0.1368421052631579
$ws.Range("F19").Value = 0.009538950715421303
$ws.Range("H19").Value = 0.1979332273449921
$ws.Range("I19").Value = 0.07074721780604133
$ws.Range("J19").Value = 0.3847376788553259
$ws.Range("K19").Value = 0.1144674085850556
$ws.Range("M19").Value = 0.01430842607313196
$ws.Range("N19").Value = 0.000794912559618442
$ws.Range("O19").Value = 0.06995230524642289
$ws.Range("S19").Value = 0.1375198728139905
